$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 ("Reason" column) had a typo - fix "Insturument" -> "Instrument".
$ws.Range("B4").Value = "Instrument"

# Row 5 keeps its existing text "Vehica rent" (re-set explicitly for clarity;
# this is a no-op in terms of displayed content).
$ws.Range("B5").Value = "Vehica rent"

# Move the active cell selection on the sheet from C5 to B5.
$ws.Range("B5").Select()
